$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column M: "SlotsPerDay" header (matching the existing bold/centered
# header style used by the other headers) with its value (10) below it.
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("M1").Value = "SlotsPerDay"
$ws.Range("M2").Value = 10

# Move the active selection to the newly added cell, mirroring the diff.
$ws.Range("M2").Select()
